$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the Illumina-style "__S<N>" suffix from the sample_id values
# in column D (rows 2-4) to produce cleaned sample names.
$ws.Range("D2").Value = "mNGplate11_sorted_A10_PATZ1-N"
$ws.Range("D3").Value = "mNGplate11_sorted_A11_KDELR3-C"
$ws.Range("D4").Value = "mNGplate11_sorted_A12_MYH9-C"
